# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.119.56'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '3.052.12'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -1.75%  '
$ws.Range("D9").Value = '3.052.50'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.154'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.81'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("E12").Value = '  -2.53%  '
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.01%  '
$ws.Range("E15").Value = '  +1.82%  '
$ws.Range("D16").Value = '3.553.59'
$ws.Range("E16").Value = '  -0.34%  '
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").Value = '63.091.76'
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = '3.049.60'
$ws.Range("E19").Value = '  -0.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '476.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.53%  '
$ws.Range("E22").Value = '  -1.45%  '
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("E24").Value = '  +2.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.21'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.34%  '
$ws.Range("E26").Value = '  -2.18%  '
$ws.Range("E27").Value = '  +7.31%  '
$ws.Range("E28").Value = '  +0.25%  '
$ws.Range("E29").Value = '  +0.92%  '
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.67%  '
$ws.Range("E34").Value = '  -2.51%  '
$ws.Range("E35").Value = '  +1.25%  '
$ws.Range("E36").Value = '  -2.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.24'
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = '  +0.69%  '
$ws.Range("E39").Value = '  -3.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.22'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.53%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.47'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '433.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.25%  '
$ws.Range("E43").Value = '  +1.17%  '
$ws.Range("E44").Value = '  +2.57%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").Value = '2.826.53'
$ws.Range("E46").Value = '  +1.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '38.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.16'
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = '  -1.68%  '
